$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 should take on the values that were in row 3 (the "MuSCs" target-cluster row),
# then row 3 (now a duplicate) is removed entirely.

$ws.Range("D2").Value = $ws.Range("D3").Value2

$ws.Range("G2").Value = 0.038688
$ws.Range("M2").Value = 0.011155
$ws.Range("N2").Value = 0.033465
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.00043156464
$ws.Range("R2").Value = 0.00388408176
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

$ws.Rows("3").Delete()
